$d = $word.ActiveDocument

# The bullet under the AIAA club-officer entry currently reads
# "Directed our AIAA Design Build Fly (DBF) RC plane design ...".
# The edit expands it to "Directed our 2024 AIAA Design Build Fly ...",
# i.e. it inserts "2024 " between "our " and "AIAA ".
$edit = $d.Content
$edit.Find.Execute("Directed our AIAA ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "Directed our 2024 AIAA ", 2) | Out-Null

# Word's autosave/normalisation pass re-merges adjacent runs that share
# identical direct formatting, which would otherwise fold the freshly
# inserted "2024 " text back into neighboring runs. Touch each piece's
# font explicitly (re-asserting the same Arial/complex-script font it
# already has) so Word keeps "Directed our ", "2024 " and "AIAA " as
# their own runs, matching the source edit.
$scope = $d.Range($edit.Start, $d.Content.End)

$run2024 = $scope.Duplicate
$run2024.Find.Execute("2024 ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$run2024.Font.Name = "Arial"
$run2024.Font.NameBi = "Arial"

$runAIAA = $scope.Duplicate
$runAIAA.Find.Execute("AIAA ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$runAIAA.Font.Name = "Arial"
$runAIAA.Font.NameBi = "Arial"
